# Applies the PtX_demand_EL.xlsx edit: insert "Fossil Gases" and "Fossil Liquids"
# rows into each year block (2030/2040/2050) and refresh the recalculated
# "Biogenic Liquids", "Synthetic Liquids" and "Overall Demand" totals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the 6 new rows from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(29).Insert()
$ws.Rows.Item(27).Insert()
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(7).Insert()

# Populate the newly inserted rows.
# Row 7: Fossil Gases (2030)
$ws.Range("A7").Value = "Fossil Gases"
$ws.Range("B7").Value = 2030
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = [double]"0.0026688655762997"
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = [double]"0.0001719976325194742"
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""

# Row 10: Fossil Liquids (2030)
$ws.Range("A10").Value = "Fossil Liquids"
$ws.Range("B10").Value = 2030
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = [double]"0.1339741348927592"
$ws.Range("G10").Value = [double]"0.0002284665558812"
$ws.Range("H10").Value = [double]"0.0987525188613676"
$ws.Range("I10").Value = [double]"0.051339213366296"
$ws.Range("J10").Value = [double]"6.910669335085698e-05"
$ws.Range("K10").Value = [double]"0.0911124665966538"

# Row 19: Fossil Gases (2040)
$ws.Range("A19").Value = "Fossil Gases"
$ws.Range("B19").Value = 2040
$ws.Range("C19").Value = ""
$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("F19").Value = [double]"0.0014641080342316"
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = ""
$ws.Range("I19").Value = [double]"0.0001817098307918636"
$ws.Range("J19").Value = ""
$ws.Range("K19").Value = ""

# Row 22: Fossil Liquids (2040)
$ws.Range("A22").Value = "Fossil Liquids"
$ws.Range("B22").Value = 2040
$ws.Range("C22").Value = ""
$ws.Range("D22").Value = ""
$ws.Range("E22").Value = ""
$ws.Range("F22").Value = [double]"0.0427780115978808"
$ws.Range("G22").Value = [double]"0.0002456142496622"
$ws.Range("H22").Value = [double]"0.09318696107376309"
$ws.Range("I22").Value = [double]"0.0227762054155748"
$ws.Range("J22").Value = [double]"6.132600771488283e-05"
$ws.Range("K22").Value = [double]"0.0881296132116407"

# Row 31: Fossil Gases (2050)
$ws.Range("A31").Value = "Fossil Gases"
$ws.Range("B31").Value = 2050
$ws.Range("C31").Value = ""
$ws.Range("D31").Value = ""
$ws.Range("E31").Value = ""
$ws.Range("F31").Value = [double]"0.0001017146949104807"
$ws.Range("G31").Value = ""
$ws.Range("H31").Value = ""
$ws.Range("I31").Value = [double]"6.895719913371703e-05"
$ws.Range("J31").Value = ""
$ws.Range("K31").Value = ""

# Row 34: Fossil Liquids (2050)
$ws.Range("A34").Value = "Fossil Liquids"
$ws.Range("B34").Value = 2050
$ws.Range("C34").Value = ""
$ws.Range("D34").Value = ""
$ws.Range("E34").Value = ""
$ws.Range("F34").Value = [double]"0.008455469550090499"
$ws.Range("G34").Value = [double]"0.000221533651345"
$ws.Range("H34").Value = [double]"0.0837043172234941"
$ws.Range("I34").Value = [double]"0.004080268444763"
$ws.Range("J34").Value = [double]"5.278641618417322e-05"
$ws.Range("K34").Value = [double]"0.0818667823546865"

# Refresh totals/values that changed because of the new Fossil Gases/Liquids rows.
# Row 9: Biogenic Liquids (2030)
$ws.Range("F9").Value = [double]"0.0136010314688469"
$ws.Range("G9").Value = [double]"3.244016311241035e-05"
$ws.Range("H9").Value = [double]"0.0108660619107528"
$ws.Range("I9").Value = [double]"0.0081773923115365"
$ws.Range("J9").Value = [double]"1.141692868523114e-05"
$ws.Range("K9").Value = [double]"0.009414850482104999"

# Row 13: Overall Demand (2030)
$ws.Range("D13").Value = [double]"0.001961873253426569"
$ws.Range("E13").Value = [double]"0.003807482572493612"
$ws.Range("F13").Value = [double]"0.151170302082553"
$ws.Range("G13").Value = [double]"0.0002609067189936103"
$ws.Range("H13").Value = [double]"0.1096185858913224"
$ws.Range("I13").Value = [double]"0.05996681229520832"
$ws.Range("J13").Value = [double]"8.052362203608813e-05"
$ws.Range("K13").Value = [double]"0.1005273170787588"

# Row 21: Biogenic Liquids (2040)
$ws.Range("F21").Value = [double]"0.0065796695125626"
$ws.Range("G21").Value = [double]"5.288281768360289e-05"
$ws.Range("H21").Value = [double]"0.0134120244731173"
$ws.Range("I21").Value = [double]"0.0053643901672684"
$ws.Range("J21").Value = [double]"1.386971540372111e-05"
$ws.Range("K21").Value = [double]"0.0107282242132517"

# Row 25: Overall Demand (2040)
$ws.Range("D25").Value = [double]"0.00210862566944133"
$ws.Range("E25").Value = [double]"0.01506913897111925"
$ws.Range("F25").Value = [double]"0.05441935047008259"
$ws.Range("G25").Value = [double]"0.0002984970673458029"
$ws.Range("H25").Value = [double]"0.1065994140806307"
$ws.Range("I25").Value = [double]"0.02873508633663617"
$ws.Range("J25").Value = [double]"7.519572311860394e-05"
$ws.Range("K25").Value = [double]"0.0988578374248924"

# Row 32: Synthetic Liquids (2050)
$ws.Range("F32").Value = [double]"1.443660391841188e-10"
$ws.Range("G32").Value = [double]"1.947388345786181e-12"
$ws.Range("H32").Value = [double]"5.016665144292069e-10"
$ws.Range("I32").Value = [double]"9.057952916775318e-11"
$ws.Range("J32").Value = [double]"9.499126337812933e-14"
$ws.Range("K32").Value = [double]"7.890151725404005e-10"

# Row 33: Biogenic Liquids (2050)
$ws.Range("F33").Value = [double]"0.001484476495781945"
$ws.Range("G33").Value = [double]"9.449231048887008e-05"
$ws.Range("H33").Value = [double]"0.0177972577647702"
$ws.Range("I33").Value = [double]"0.0013799432361925"
$ws.Range("J33").Value = [double]"1.783865034937578e-05"
$ws.Range("K33").Value = [double]"0.0152473298017653"

# Row 37: Overall Demand (2050)
$ws.Range("D37").Value = [double]"0.002229631621553799"
$ws.Range("E37").Value = [double]"0.03904990144600101"
$ws.Range("F37").Value = [double]"0.01469872163549077"
$ws.Range("G37").Value = [double]"0.0003160259637812584"
$ws.Range("H37").Value = [double]"0.1015023018201198"
$ws.Range("I37").Value = [double]"0.006092883145501675"
$ws.Range("J37").Value = [double]"7.062506662854026e-05"
$ws.Range("K37").Value = [double]"0.09711411294546697"
